# Update numbers per "Sediki di update oleh Samuel"
$wb = $excel.ActiveWorkbook

# --- Sheet: Distribusi Log Training ---
$ws1 = $wb.Worksheets.Item("Distribusi Log Training")
$ws1.Range("B2").Value = 39
$ws1.Range("C2").Value = 78
$ws1.Range("B3").Value = 6
$ws1.Range("C3").Value = 12

# --- Sheet: Metrik Akurasi Training ---
$ws2 = $wb.Worksheets.Item("Metrik Akurasi Training")
$ws2.Range("B2").Value = 39
$ws2.Range("C2").Value = 39
$ws2.Range("B3").Value = 6
$ws2.Range("C3").Value = 6

# --- Sheet: Distribusi Log Testing ---
$ws3 = $wb.Worksheets.Item("Distribusi Log Testing")
$ws3.Range("B2").Value = 16
$ws3.Range("C2").Value = 69.56521739130434
$ws3.Range("A3").Value = "Brute Force"
$ws3.Range("B3").Value = 7
$ws3.Range("C3").Value = 30.43478260869566

# --- Sheet: Metrik Akurasi Testing ---
$ws4 = $wb.Worksheets.Item("Metrik Akurasi Testing")
$ws4.Range("B2").Value = 16
$ws4.Range("C2").Value = 16
$ws4.Range("A3").Value = "Brute Force"
$ws4.Range("B3").Value = 7
$ws4.Range("C3").Value = 5
$ws4.Range("D3").Value = 71.42857142857143

# --- Sheet: Analisis IP Penyerang ---
$ws6 = $wb.Worksheets.Item("Analisis IP Penyerang")
$ws6.Range("A4").Value = "96.0.4664.110"
$ws6.Range("A5").Value = "192.168.0.10"
$ws6.Range("B5").Value = 1
$ws6.Range("C5").Value = "SQL Injection: 1"
